$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Settings"
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")

# logF_BusinessProcessName value: process renamed from "Multiple Clients" to "By State"
$wsSettings.Range("B5").Value2 = "EDI Generation - By State"

# The "NameEdiXmlFile" argument/setting row was removed entirely (it is no
# longer used because the bot won't generate that file name anymore).
$wsSettings.Rows.Item(9).Delete()

# A new setting row was added: StatusTable_LegalEntityColumn / "Legal Entity"
# It is inserted right after "StatusTable_ClientColumn" (now row 17) and
# before "StatusTable_ReturnColumn" (now row 18 prior to the insert).
$wsSettings.Rows.Item(18).Insert()
$wsSettings.Rows.Item(18).RowHeight = 14.25
$wsSettings.Range("A18").Value2 = "StatusTable_LegalEntityColumn"
$wsSettings.Range("B18").Value2 = "Legal Entity"

# Update the view state for this sheet (scroll position + selection)
$wsSettings.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$wsSettings.Range("B7").Select()

# ---------------------------------------------------------------------------
# Sheet "Constants"
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")

# Mail_Subject value: process renamed from "Multiple Clients" to "By State"
$wsConstants.Range("B37").Value2 = "EDI Generation - By State"

# Mail_HtmlBody_Header value: bot no longer zips/attaches EDI files, instead
# it points to the output folder, and it no longer deletes that folder.
$wsConstants.Range("B38").Value2 = 'Hi team, <br><br>Here is the status of the bot running "EDI Generation - By State" process. <br>You can find the Output folder here: <b>C:\Users\<USERNAME>\Desktop\EDI Generation output<b> <br>Here is the status table:'

# Update the view state for this sheet (scroll position + selection)
$wsConstants.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$wsConstants.Range("C38").Select()
